$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 22.700661
$ws.Cells.Item(2, 8).Value = 68.10198299999999
$ws.Cells.Item(2, 9).Value = 0.08615268874617349
$ws.Cells.Item(2, 10).Value = 0.08615268874617349
$ws.Cells.Item(2, 13).Value = 0.967553
$ws.Cells.Item(2, 14).Value = 2.902659
$ws.Cells.Item(2, 15).Value = 0.01359591950841534
$ws.Cells.Item(2, 16).Value = 0.01359591950841534
$ws.Cells.Item(2, 17).Value = 21.964092652533
$ws.Cells.Item(2, 18).Value = 197.676833872797
$ws.Cells.Item(2, 19).Value = 0.001171325021626535
$ws.Cells.Item(2, 20).Value = 0.001171325021626535
$ws.Cells.Item(3, 7).Value = 22.700661
$ws.Cells.Item(3, 8).Value = 68.10198299999999
$ws.Cells.Item(3, 9).Value = 0.08615268874617349
$ws.Cells.Item(3, 10).Value = 0.08615268874617349
$ws.Cells.Item(3, 13).Value = 5.176377666666667
$ws.Cells.Item(3, 15).Value = 0.07273773540173906
$ws.Cells.Item(3, 16).Value = 0.07273773540173906
$ws.Cells.Item(3, 17).Value = 117.507194618971
$ws.Cells.Item(3, 18).Value = 1057.564751570739
$ws.Cells.Item(3, 19).Value = 0.006266551478167549
$ws.Cells.Item(3, 20).Value = 0.006266551478167549
$ws.Cells.Item(4, 7).Value = 22.700661
$ws.Cells.Item(4, 8).Value = 68.10198299999999
$ws.Cells.Item(4, 9).Value = 0.08615268874617349
$ws.Cells.Item(4, 10).Value = 0.08615268874617349
$ws.Cells.Item(4, 13).Value = 4.670153
$ws.Cells.Item(4, 14).Value = 14.010459
$ws.Cells.Item(4, 15).Value = 0.06562433714740633
$ws.Cells.Item(4, 16).Value = 0.06562433714740633
$ws.Cells.Item(4, 17).Value = 106.015560071133
$ws.Cells.Item(4, 18).Value = 954.1400406401968
$ws.Cells.Item(4, 19).Value = 0.005653713092434448
$ws.Cells.Item(4, 20).Value = 0.005653713092434448
$ws.Cells.Item(5, 7).Value = 22.700661
$ws.Cells.Item(5, 8).Value = 68.10198299999999
$ws.Cells.Item(5, 9).Value = 0.08615268874617349
$ws.Cells.Item(5, 10).Value = 0.08615268874617349
$ws.Cells.Item(5, 13).Value = 60.35087133333334
$ws.Cells.Item(5, 14).Value = 181.052614
$ws.Cells.Item(5, 15).Value = 0.8480420079424392
$ws.Cells.Item(5, 16).Value = 0.8480420079424392
$ws.Cells.Item(5, 17).Value = 1370.004671192618
$ws.Cells.Item(5, 18).Value = 12330.04204073356
$ws.Cells.Item(5, 19).Value = 0.07306109915394496
$ws.Cells.Item(5, 20).Value = 0.07306109915394496
$ws.Cells.Item(6, 9).Value = 0.5030288587986086
$ws.Cells.Item(6, 10).Value = 0.5030288587986087
$ws.Cells.Item(6, 13).Value = 0.967553
$ws.Cells.Item(6, 14).Value = 2.902659
$ws.Cells.Item(6, 15).Value = 0.01359591950841534
$ws.Cells.Item(6, 16).Value = 0.01359591950841534
$ws.Cells.Item(6, 17).Value = 128.2440817848683
$ws.Cells.Item(6, 18).Value = 1154.196736063815
$ws.Cells.Item(6, 19).Value = 0.006839139874635909
$ws.Cells.Item(6, 20).Value = 0.006839139874635911
$ws.Cells.Item(7, 9).Value = 0.5030288587986086
$ws.Cells.Item(7, 10).Value = 0.5030288587986087
$ws.Cells.Item(7, 13).Value = 5.176377666666667
$ws.Cells.Item(7, 15).Value = 0.07273773540173906
$ws.Cells.Item(7, 16).Value = 0.07273773540173906
$ws.Cells.Item(7, 17).Value = 686.1017441249895
$ws.Cells.Item(7, 19).Value = 0.03658918003073195
$ws.Cells.Item(7, 20).Value = 0.03658918003073196
$ws.Cells.Item(8, 9).Value = 0.5030288587986086
$ws.Cells.Item(8, 10).Value = 0.5030288587986087
$ws.Cells.Item(8, 13).Value = 4.670153
$ws.Cells.Item(8, 14).Value = 14.010459
$ws.Cells.Item(8, 15).Value = 0.06562433714740633
$ws.Cells.Item(8, 16).Value = 0.06562433714740633
$ws.Cells.Item(8, 17).Value = 619.0043163318683
$ws.Cells.Item(8, 18).Value = 5571.038846986814
$ws.Cells.Item(8, 19).Value = 0.03301093542467494
$ws.Cells.Item(8, 20).Value = 0.03301093542467495
$ws.Cells.Item(9, 9).Value = 0.5030288587986086
$ws.Cells.Item(9, 10).Value = 0.5030288587986087
$ws.Cells.Item(9, 13).Value = 60.35087133333334
$ws.Cells.Item(9, 14).Value = 181.052614
$ws.Cells.Item(9, 15).Value = 0.8480420079424392
$ws.Cells.Item(9, 16).Value = 0.8480420079424392
$ws.Cells.Item(9, 17).Value = 7999.191857252333
$ws.Cells.Item(9, 18).Value = 71992.72671527098
$ws.Cells.Item(9, 19).Value = 0.4265896034685658
$ws.Cells.Item(9, 20).Value = 0.4265896034685658
$ws.Cells.Item(10, 7).Value = 41.94534433333333
$ws.Cells.Item(10, 8).Value = 125.836033
$ws.Cells.Item(10, 9).Value = 0.159189381961201
$ws.Cells.Item(10, 10).Value = 0.159189381961201
$ws.Cells.Item(10, 13).Value = 0.967553
$ws.Cells.Item(10, 14).Value = 2.902659
$ws.Cells.Item(10, 15).Value = 0.01359591950841534
$ws.Cells.Item(10, 16).Value = 0.01359591950841534
$ws.Cells.Item(10, 17).Value = 40.58434374574966
$ws.Cells.Item(10, 18).Value = 365.259093711747
$ws.Cells.Item(10, 19).Value = 0.002164326023738874
$ws.Cells.Item(10, 20).Value = 0.002164326023738874
$ws.Cells.Item(11, 7).Value = 41.94534433333333
$ws.Cells.Item(11, 8).Value = 125.836033
$ws.Cells.Item(11, 9).Value = 0.159189381961201
$ws.Cells.Item(11, 10).Value = 0.159189381961201
$ws.Cells.Item(11, 13).Value = 5.176377666666667
$ws.Cells.Item(11, 15).Value = 0.07273773540173906
$ws.Cells.Item(11, 16).Value = 0.07273773540173906
$ws.Cells.Item(11, 17).Value = 217.1249436277099
$ws.Cells.Item(11, 18).Value = 1954.124492649389
$ws.Cells.Item(11, 19).Value = 0.01157907514386021
$ws.Cells.Item(11, 20).Value = 0.01157907514386021
$ws.Cells.Item(12, 7).Value = 41.94534433333333
$ws.Cells.Item(12, 8).Value = 125.836033
$ws.Cells.Item(12, 9).Value = 0.159189381961201
$ws.Cells.Item(12, 10).Value = 0.159189381961201
$ws.Cells.Item(12, 13).Value = 4.670153
$ws.Cells.Item(12, 14).Value = 14.010459
$ws.Cells.Item(12, 15).Value = 0.06562433714740633
$ws.Cells.Item(12, 16).Value = 0.06562433714740633
$ws.Cells.Item(12, 17).Value = 195.8911756743497
$ws.Cells.Item(12, 18).Value = 1763.020581069147
$ws.Cells.Item(12, 19).Value = 0.0104466976721091
$ws.Cells.Item(12, 20).Value = 0.0104466976721091
$ws.Cells.Item(13, 7).Value = 41.94534433333333
$ws.Cells.Item(13, 8).Value = 125.836033
$ws.Cells.Item(13, 9).Value = 0.159189381961201
$ws.Cells.Item(13, 10).Value = 0.159189381961201
$ws.Cells.Item(13, 13).Value = 60.35087133333334
$ws.Cells.Item(13, 14).Value = 181.052614
$ws.Cells.Item(13, 15).Value = 0.8480420079424392
$ws.Cells.Item(13, 16).Value = 0.8480420079424392
$ws.Cells.Item(13, 17).Value = 2531.438078893363
$ws.Cells.Item(13, 18).Value = 22782.94271004026
$ws.Cells.Item(13, 19).Value = 0.1349992831214928
$ws.Cells.Item(13, 20).Value = 0.1349992831214928
$ws.Cells.Item(14, 7).Value = 66.302588
$ws.Cells.Item(14, 8).Value = 198.907764
$ws.Cells.Item(14, 9).Value = 0.2516290704940168
$ws.Cells.Item(14, 10).Value = 0.2516290704940168
$ws.Cells.Item(14, 13).Value = 0.967553
$ws.Cells.Item(14, 14).Value = 2.902659
$ws.Cells.Item(14, 15).Value = 0.01359591950841534
$ws.Cells.Item(14, 16).Value = 0.01359591950841534
$ws.Cells.Item(14, 17).Value = 64.151267927164
$ws.Cells.Item(14, 18).Value = 577.361411344476
$ws.Cells.Item(14, 19).Value = 0.003421128588414023
$ws.Cells.Item(14, 20).Value = 0.003421128588414023
$ws.Cells.Item(15, 7).Value = 66.302588
$ws.Cells.Item(15, 8).Value = 198.907764
$ws.Cells.Item(15, 9).Value = 0.2516290704940168
$ws.Cells.Item(15, 10).Value = 0.2516290704940168
$ws.Cells.Item(15, 13).Value = 5.176377666666667
$ws.Cells.Item(15, 15).Value = 0.07273773540173906
$ws.Cells.Item(15, 16).Value = 0.07273773540173906
$ws.Cells.Item(15, 17).Value = 343.2072357654014
$ws.Cells.Item(15, 18).Value = 3088.865121888612
$ws.Cells.Item(15, 19).Value = 0.01830292874897934
$ws.Cells.Item(15, 20).Value = 0.01830292874897934
$ws.Cells.Item(16, 7).Value = 66.302588
$ws.Cells.Item(16, 8).Value = 198.907764
$ws.Cells.Item(16, 9).Value = 0.2516290704940168
$ws.Cells.Item(16, 10).Value = 0.2516290704940168
$ws.Cells.Item(16, 13).Value = 4.670153
$ws.Cells.Item(16, 14).Value = 14.010459
$ws.Cells.Item(16, 15).Value = 0.06562433714740633
$ws.Cells.Item(16, 16).Value = 0.06562433714740633
$ws.Cells.Item(16, 17).Value = 309.643230255964
$ws.Cells.Item(16, 18).Value = 2786.789072303676
$ws.Cells.Item(16, 19).Value = 0.01651299095818784
$ws.Cells.Item(16, 20).Value = 0.01651299095818784
$ws.Cells.Item(17, 7).Value = 66.302588
$ws.Cells.Item(17, 8).Value = 198.907764
$ws.Cells.Item(17, 9).Value = 0.2516290704940168
$ws.Cells.Item(17, 10).Value = 0.2516290704940168
$ws.Cells.Item(17, 13).Value = 60.35087133333334
$ws.Cells.Item(17, 14).Value = 181.052614
$ws.Cells.Item(17, 15).Value = 0.8480420079424392
$ws.Cells.Item(17, 16).Value = 0.8480420079424392
$ws.Cells.Item(17, 17).Value = 4001.418957455011
$ws.Cells.Item(17, 18).Value = 36012.77061709509
$ws.Cells.Item(17, 19).Value = 0.2133920221984356
$ws.Cells.Item(17, 20).Value = 0.2133920221984356

$wb.Save()